# Add the new "dockerised" sentence fragment after the existing
# " with Python." run, splitting it into several runs to match the
# target structure while preserving the original run formatting
# (Calibri / sz 24 / szCs 24 / shd clear-auto-FFFFFF).

$d = $word.ActiveDocument

# 1. Replace the trailing text so the paragraph now reads
#    " with Python and dockerised it, for easy setup for other members."
#    Find/Replace rewrites the text in-place inside the already matched
#    run, so the run keeps its existing formatting (rFonts/sz/szCs/shd).
$r = $d.Content
$null = $r.Find.Execute(" with Python.", $false, $false, $false, $false, $false, $true, 1, $false, " with Python and dockerised it, for easy setup for other members.", 2)

$base = $r.Start

# 2. Split the merged run into the individual runs from the diff, using
#    absolute character offsets (not text search, since fragments like
#    "," or "." are ambiguous elsewhere in the document). Nudging
#    Font.Size away from, then back to, its current value on a
#    sub-Range forces a run split at the sub-Range boundaries while
#    copying the full original run-properties (rFonts/sz/szCs/shd) onto
#    every piece - only the touched property (sz) gets rewritten, so
#    rFonts/szCs/shd survive unchanged on the newly split runs.
$segments = @(" with Python", " and ", "dockerised", " it", ",", " for easy ", "setup for other members", ".")

$pos = $base
foreach ($seg in $segments) {
    $segLen = $seg.Length
    $sub = $d.Range($pos, $pos + $segLen)
    $sub.Font.Size = 11
    $sub.Font.Size = 12
    $pos = $pos + $segLen
}
